$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update Total row marks/score values
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 45
$ws.Range("E12").Value = "45/140"
